$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Rename the four inline pictures (docPr/name swap).
#    InlineShape has no writable "Name" in the Word object model;
#    only Shape does, so we round-trip through ConvertToShape /
#    ConvertToInlineShape to rename, then restore it as an inline
#    picture exactly as it was.
# ---------------------------------------------------------------
$renames = @{
    4 = "image2.png"
    5 = "image4.png"
    2 = "image1.png"
    6 = "image3.png"
}

$count = $d.InlineShapes.Count
for ($i = 1; $i -le $count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shape = $shp.ConvertToShape()
    $id = $shape.ID
    if ($renames.ContainsKey($id)) {
        $shape.Name = $renames[$id]
    }
    [void]$shape.ConvertToInlineShape()
}

# ---------------------------------------------------------------
# 2. Update the SQL snippet to use CONCAT() and drop the quoting
#    around the typeSearch parameter.
# ---------------------------------------------------------------
$newText = '    `name` LIKE CONCAT(''%'', nameSearch, ''%'') AND `type` = typeSearch'

$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*nameSearch*") {
        $prng = $p.Range
        $prng.MoveEnd(1, -1)
        $prng.Text = $newText
        break
    }
}
